$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 36
$ws.Range("F4").Value = 843
$ws.Range("F7").Value = 9620
$ws.Range("F8").Value = 40
$ws.Range("F9").Value = 719
$ws.Range("F10").Value = 2181
$ws.Range("F12").Value = 1651
$ws.Range("F13").Value = 2777
$ws.Range("F14").Value = 145
$ws.Range("F15").Value = 4146
$ws.Range("F16").Value = 350
$ws.Range("F17").Value = 173
$ws.Range("F18").Value = 136
$ws.Range("F19").Value = 526
$ws.Range("F20").Value = 245
$ws.Range("F21").Value = 38
$ws.Range("F23").Value = 87
$ws.Range("F24").Value = 284
$ws.Range("F25").Value = 4009
$ws.Range("F27").Value = 3479
$ws.Range("F28").Value = 1108
$ws.Range("F29").Value = 204
$ws.Range("F30").Value = 507
$ws.Range("F31").Value = 4350
$ws.Range("F32").Value = 76
$ws.Range("F33").Value = 343
$ws.Range("F34").Value = 432
$ws.Range("F35").Value = 327

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1010

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 36
$ws.Range("F4").Value = 1010
$ws.Range("F6").Value = 843
$ws.Range("F9").Value = 9620
$ws.Range("F10").Value = 40
$ws.Range("F11").Value = 719
$ws.Range("F12").Value = 2182
$ws.Range("F14").Value = 1651
$ws.Range("F16").Value = 2777
$ws.Range("F17").Value = 145
$ws.Range("F18").Value = 4146
$ws.Range("F19").Value = 350
$ws.Range("F20").Value = 173
$ws.Range("F21").Value = 136
$ws.Range("F22").Value = 526
$ws.Range("F23").Value = 245
$ws.Range("F24").Value = 38
$ws.Range("F27").Value = 87
$ws.Range("F28").Value = 284
$ws.Range("F29").Value = 4009
$ws.Range("F31").Value = 3479
$ws.Range("F32").Value = 1108
$ws.Range("F33").Value = 204
$ws.Range("F34").Value = 507
$ws.Range("F35").Value = 4350
$ws.Range("F36").Value = 76
$ws.Range("F37").Value = 343
$ws.Range("F38").Value = 432
$ws.Range("F39").Value = 327
$ws.Range("F41").Value = 25
